# Budget Summary update: fill in real event-budget data and expand the
# expense breakdown, matching the "added budget and chatbot controller" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells whose old content is being relocated / removed ---------
$ws.Range("A6").ClearContents()
$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("A11").ClearContents()

# --- Header block (rows 1-4) --------------------------------------------
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Amount"

$ws.Range("A2").Value = "Your name"
$ws.Range("B2").Value = "John Doe"

$ws.Range("A3").Value = "Booking Artist"
$ws.Range("B3").Value = "Wizkid"

$ws.Range("A4").Value = "Location"
$ws.Range("B4").Value = "London UK"

# --- Capacity / date / currency (rows 5-7) ------------------------------
$ws.Range("A5").Value = "Capacity"
$ws.Range("B5").Value = 1000

# B6 must stay literal text "2025-06-19" (not be auto-converted to a date
# serial number). Build it as a formula that evaluates to that string and
# then convert the cell to a plain value, which keeps the cell's existing
# (non-date) number format/style instead of Excel stamping a date format
# on it.
$ws.Range("B6").Formula = '="2025-06-19"'
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4163) | Out-Null

$ws.Range("A7").Value = "Currency"
$ws.Range("B7").Value = "USD"

# --- Ticket tiers (rows 9-10) --------------------------------------------
$ws.Range("A9").Value = "Ticket Tiers"
$ws.Range("B9").Value = 1

$ws.Range("A10").Value = "Regular (1000 @ 100)"
$ws.Range("B10").Value = 100000

# --- Expense breakdown (rows 12-32), all 5000 each -----------------------
$expenseItems = @(
    "Venue",
    "Venue barricade",
    "Box office Personnel",
    "Ticket Printing",
    "Ushers",
    "Security",
    "Lighting Technician",
    "Audio Technician",
    "Stage Hands",
    "Medical Cost",
    "Stage",
    "Lights?",
    "Backline",
    "Sound",
    "Outdoor Posters",
    "Radio Advertising",
    "Artist Fee",
    "Digital Ads",
    "Flights",
    "Accommodation",
    "Catering"
)

$row = 12
foreach ($item in $expenseItems) {
    $ws.Cells.Item($row, 1).Value = $item
    $ws.Cells.Item($row, 2).Value = 5000
    $row = $row + 1
}

# --- Totals (rows 34-36) ---------------------------------------------------
$ws.Range("A34").Value = "Gross Revenue"
$ws.Range("B34").Value = 100000

$ws.Range("A35").Value = "Total Expenses"
$ws.Range("B35").Value = 105000

$ws.Range("A36").Value = "Net Profit"
$ws.Range("B36").Value = -5000

Write-Output "Budget summary populated"
